$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("F3").Value = "U3-110"
$ws.Range("F6").Value = "U3-110"
$ws.Range("F9").Value = "U3-110"
$ws.Range("F11").Value = "U4-307"
$ws.Range("F12").Value = "U4-308"
$ws.Range("F13").Value = "U4-307"
$ws.Range("F14").Value = "U4-308"
$ws.Range("F17").Value = "U3-Amphi"
$ws.Range("F18").Value = "U3-Amphi"
$ws.Range("F19").Value = "U1-Mathis"
$ws.Range("F20").Value = "U1-Mathis"

$wb.Save()
